$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Document TileMap" task is finished; C2 now shows what used to be C3
# ("Document SDLUtilityTool"), and the old C3 row is vacated.
$ws.Range("C2").Value = "Document SDLUtilityTool"
$ws.Range("C3").Value = $null

# New agenda item discovered while documenting the TileMap.
$ws.Range("A4").Value = "Figure out level file format"

# Minimum tile world size determination got reprioritized -- mark it.
$ws.Range("A3").Value = "*Determine Texture Sizes"

# Keep the active selection tidy (matches the single C2 selection in workbook).
$ws.Range("C2").Select()
